$d = $word.ActiveDocument

# Locate the full sentence to edit: ". In many cases there is not single ...decisions."
# It currently spans 3 runs, with "not" wrapped in proofErr spellStart/spellEnd
# tags (so its two boundaries are "protected"); the other run boundaries are
# plain <w:r> adjacency.
$anchor = $d.Content
$oldSentence = ". In many cases there is not single " + [char]8220 + "correct" + [char]8221 + " answer. Instead, be prepared to explain your decisions."
$found = $anchor.Find.Execute($oldSentence, $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate target sentence in document."
}
$start = $anchor.Start
$end = $anchor.End

$newSentence = ". In many cases there is not a single " + [char]8220 + "correct" + [char]8221 + " answer. Instead, be prepared to explain your decisions."

# Replace the whole matched span in one go. Because the matched range
# crosses both proofErr boundaries around "not", this also removes the
# (now stale) proofErr spellStart/spellEnd markers, and merges everything
# touched (including the preceding run) into a single run.
$full = $d.Range($start, $end)
$full.Text = $newSentence

$newEnd = $start + $newSentence.Length

# Re-split the merged run back into separate runs at the boundaries we
# actually want, using a harmless formatting round-trip (bold on, then
# off) which forces a clean run split without merging text across it
# (unlike plain .Text assignment, which merges all touched runs).
function Split-At($pos) {
    $r = $d.Range($pos, $newEnd)
    $r.Font.Bold = 1
    $r.Font.Bold = 0
}

# boundary restoring the original run break right before ". In many cases..."
Split-At $start
# boundary right after "...there is not " (before "a ")
$boundary1 = $start + ". In many cases there is not ".Length
Split-At $boundary1
# boundary right after "a " (before "single...")
$boundary2 = $boundary1 + "a ".Length
Split-At $boundary2
